$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.007307977243890207
$ws.Range("D2").Value = 0.2042482159442283
$ws.Range("E2").Value = 4.734106030302786
$ws.Range("F2").Value = [double]"1.737961192307855e-05"
$ws.Range("G2").Value = 0.9987761858899808
$ws.Range("H2").Value = 0.380665828052104
$ws.Range("I2").Value = 5.825676796000084
$ws.Range("J2").Value = [double]"1.516402179975824e-15"
$ws.Range("K2").Value = 41
$ws.Range("B3").Value = 0.9999999999999998
$ws.Range("C3").Value = 0.1550479081649281
$ws.Range("D3").Value = 0.1161157312238916
$ws.Range("E3").Value = 2.504872879457793
$ws.Range("F3").Value = [double]"2.336198482367304e-09"
$ws.Range("G3").Value = 0.9510833923457575
$ws.Range("H3").Value = 0.7259421860880491
$ws.Range("I3").Value = 3.058852846571612
$ws.Range("J3").Value = [double]"4.736951571734001e-15"
$ws.Range("K3").Value = 30
$ws.Range("C4").Value = 7.548961188007898
$ws.Range("D4").Value = 0.155068947513865
$ws.Range("E4").Value = 3.58098173330632
$ws.Range("F4").Value = [double]"1.554158763801489e-07"
$ws.Range("G4").Value = 0.0418586844402117
$ws.Range("H4").Value = 0.529180045017565
$ws.Range("I4").Value = 2.754992322690183
$ws.Range("J4").Value = [double]"-5.192427684400732e-15"
$ws.Range("K4").Value = 39
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 3.991571880794022
$ws.Range("D5").Value = 0.146500487672622
$ws.Range("E5").Value = 3.47558585552411
$ws.Range("F5").Value = [double]"2.045380335021598e-07"
$ws.Range("G5").Value = 0.2604958087660011
$ws.Range("H5").Value = 0.6246302966024428
$ws.Range("I5").Value = 3.217232981367922
$ws.Range("J5").Value = [double]"-1.474376176702208e-14"
$ws.Range("K5").Value = 30
$ws.Range("C6").Value = 6.380809795566059
$ws.Range("D6").Value = 0.1435430695320799
$ws.Range("E6").Value = 3.265438079660177
$ws.Range("F6").Value = [double]"1.418164207188511e-07"
$ws.Range("G6").Value = 0.06075071653480218
$ws.Range("H6").Value = 0.6341438957377079
$ws.Range("I6").Value = 2.909453587146678
$ws.Range("J6").Value = [double]"-1.01844458792281e-14"
$ws.Range("K6").Value = 30
$ws.Range("C7").Value = 3.567888816915884
$ws.Range("D7").Value = 0.1433839153537295
$ws.Range("E7").Value = 3.366916192389526
$ws.Range("F7").Value = [double]"4.123936309411831e-08"
$ws.Range("G7").Value = 0.2965439946159347
$ws.Range("H7").Value = 0.5815430501635944
$ws.Range("I7").Value = 3.786860095324411
$ws.Range("J7").Value = [double]"2.592520792638203e-15"
$ws.Range("K7").Value = 37
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.02884030168999363
$ws.Range("D8").Value = 0.1965765756382262
$ws.Range("E8").Value = 4.585481204082646
$ws.Range("F8").Value = [double]"1.673343332618078e-05"
$ws.Range("G8").Value = 0.9950220535642406
$ws.Range("H8").Value = 0.4549770928067196
$ws.Range("I8").Value = 5.221496034255939
$ws.Range("J8").Value = [double]"7.105427357601002e-15"
$ws.Range("K8").Value = 33
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 5.086357067735555
$ws.Range("D9").Value = 0.06930718222530308
$ws.Range("E9").Value = 1.480981933674747
$ws.Range("F9").Value = [double]"9.108925829895073e-15"
$ws.Range("G9").Value = 0.0018102221604539
$ws.Range("H9").Value = 0.8777311088258974
$ws.Range("I9").Value = 1.723169989892944
$ws.Range("J9").Value = [double]"-1.260640337638887e-14"
$ws.Range("K9").Value = 31
